$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 70
$ws.Range("E2").Value = 70
$ws.Range("G2").Value = 70
$ws.Range("I2").Value = 70
$ws.Range("J2").Value = 70

$ws.Range("J2").Select()
